$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing Status cells (J column) ---
# Rows 2-5: Open -> Done
$ws.Range("J2").Value = "Done"
$ws.Range("J3").Value = "Done"
$ws.Range("J4").Value = "Done"
$ws.Range("J5").Value = "Done"

# --- Add the new rows (53-61) in order so new shared strings line up ---
$ws.Range("A53").Value = 52
$ws.Range("B53").Value = "Delete table when post deleted from admin page"
$ws.Range("C53").Value = "Infrastructure"
$ws.Range("E53").Value = "Medium"
$ws.Range("F53").Value = "Low"
$ws.Range("J53").Value = "Open"

$ws.Range("A54").Value = 53
$ws.Range("B54").Value = "Job to remove orphan tables"
$ws.Range("C54").Value = "Infrastructure"
$ws.Range("E54").Value = "Summary"
$ws.Range("F54").Value = "Low"
$ws.Range("J54").Value = "Open"

$ws.Range("A55").Value = 54
$ws.Range("B55").Value = "Create admin page to maintain tables"
$ws.Range("C55").Value = "Feature"
$ws.Range("E55").Value = "Summary"
$ws.Range("F55").Value = "Low"
$ws.Range("G55").Value = 1.1
$ws.Range("J55").Value = "Open"

$ws.Range("A56").Value = 55
$ws.Range("B56").Value = "Create table from CSV file"
$ws.Range("C56").Value = "Feature"
$ws.Range("E56").Value = "Summary"
$ws.Range("F56").Value = "Low"
$ws.Range("G56").Value = 1.1
$ws.Range("J56").Value = "Open"

$ws.Range("A57").Value = 56
$ws.Range("B57").Value = "Determine if useEffect for # col/row is still needed"
$ws.Range("C57").Value = "QA"
$ws.Range("D57").Value = "App"
$ws.Range("E57").Value = "Detailed"
$ws.Range("F57").Value = "High"
$ws.Range("G57").Value = 1
$ws.Range("J57").Value = "Open"

$ws.Range("A58").Value = 57
$ws.Range("B58").Value = "update edit.js to remove depricated useSetting"
$ws.Range("C58").Value = "QA"
$ws.Range("D58").Value = "App"
$ws.Range("E58").Value = "Detailed"
$ws.Range("F58").Value = "High"
$ws.Range("G58").Value = 1
$ws.Range("J58").Value = "Open"

$ws.Range("A59").Value = 58
$ws.Range("B59").Value = "Remove red table border from block editor"
$ws.Range("C59").Value = "QA"
$ws.Range("D59").Value = "App"
$ws.Range("E59").Value = "Detailed"
$ws.Range("F59").Value = "High"
$ws.Range("G59").Value = 1
$ws.Range("J59").Value = "Open"

$ws.Range("A60").Value = 59
$ws.Range("B60").Value = 'Gray "Freeze  Header Row" when Headers are disabled'
$ws.Range("C60").Value = "QA"
$ws.Range("D60").Value = "App"
$ws.Range("E60").Value = "Detailed"
$ws.Range("F60").Value = "High"
$ws.Range("G60").Value = 1
$ws.Range("J60").Value = "Open"

$ws.Range("B61").Value = "Remove the Site Save Message or make it meaningful"

# J7: Open -> Testing (new shared string, added last so it lands at the end)
$ws.Range("J7").Value = "Testing"

# --- Column B width (widened to fit the new, longer text) ---
$ws.Columns.Item(2).ColumnWidth = 58.83

# --- Row 2 becomes hidden (status no longer "Open") ---
$ws.Rows.Item(2).Hidden = $true

# --- Update the selected/active cell ---
$ws.Range("K8").Select() | Out-Null
